# "casos de uso web"
# Adds the new "Servicios"/"Eventos"/"prototipo" catalogue rows to the
# "Elementos del Proyecto" sheet and relabels a couple of existing
# "usuario(s)" use-case rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# 1) Stamp the F:H formatting used by the existing "Casos de Uso" rows
#    (border + wrap on F, border only on G/H) onto the new rows 13-24
#    before filling in their values, so the new cells pick up the same
#    cell styles already present in the workbook (no new style slots).
# ------------------------------------------------------------------
$ws.Range("F11:H11").Copy()
$ws.Range("F13:H24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Rows 13-15: Servicios (alta/listar/buscar)
$ws.Range("F13").Value = "Alta servicio"
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = 6

$ws.Range("F14").Value = "Listar Servicios"
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = 6

$ws.Range("F15").Value = "Buscar Servicios"
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 6

# (entered out of row order - editar before borrar - to match the
# original authoring session)
$ws.Range("F17").Value = "editar Servicios"
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 6

$ws.Range("F16").Value = "borrar servicios"
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 6

# ------------------------------------------------------------------
# 2) Relabel the existing usuarios rows (12 then 11), reusing the
#    plain "no border / no wrap" style already used elsewhere in the
#    workbook (Trazabilidad sheet F2) for row 12.
# ------------------------------------------------------------------
$ws2.Range("F2").Copy()
$ws.Range("F12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F12").Value = "Editar usuarios"

$ws.Range("F11").Value = "borrar usuarios"

# Rows 18-24: Eventos, logout, estadisticas
$ws.Range("F18").Value = "Alta evento"
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 7

$ws.Range("F19").Value = "Listar eventos"
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 7

$ws.Range("F20").Value = "Buscar eventos"
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 7

$ws.Range("F21").Value = "borrar eventos"
$ws.Range("G21").Value = 18
$ws.Range("H21").Value = 7

$ws.Range("F22").Value = "editar eventos"
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = 7

$ws.Range("F23").Value = "logout"
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 8

$ws.Range("F24").Value = "estadisticas"
$ws.Range("G24").Value = 21
$ws.Range("H24").Value = 9

# ------------------------------------------------------------------
# 3) Fill in the "Interfaz de Usuario" block (J4:L9) - the prototype
#    names, their IDs and the use-case-id lists they cover.
# ------------------------------------------------------------------
$ws.Range("J4").Value = "prototipo usuarios"
$ws.Range("L4").Value = "5,6,7,8,9"
$ws.Range("J5").Value = "prototipo login"
$ws.Range("J6").Value = "prototipo eventos"
$ws.Range("J7").Value = "prototipo servicios"
$ws.Range("J8").Value = "prototipo web publica"
$ws.Range("J9").Value = "prototipo estadisticas"
$ws.Range("L6").Value = "15,16,17,19,19"
$ws.Range("L7").Value = "10,11,12,13,14"
$ws.Range("L8").Value = "1,2,3"

$ws.Range("K4").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 4
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 4
$ws.Range("K8").Value = 5
$ws.Range("K9").Value = 6
$ws.Range("L9").Value = 21

# ------------------------------------------------------------------
# 4) Restore the view: scrolled so column E is left-most, L10 selected.
# ------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("L10").Select()
